# Auto-generated COM-interop script applying the cryptos.xlsx price/volume refresh
# (commit: 'Updated cryptos list on Wed Jun 12 14:38:03 UTC 2024 with GitHub Actions').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.942.23'
$ws.Cells.Item(2, 5).Value = '  +4.54%  '
$ws.Cells.Item(3, 4).Value = '3.638.09'
$ws.Cells.Item(3, 5).Value = '  +3.36%  '
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).Value = '''634.51'
$ws.Cells.Item(5, 5).Value = '  +4.13%  '
$ws.Cells.Item(6, 4).Value = '''160.47'
$ws.Cells.Item(6, 5).Value = '  +5.61%  '
$ws.Cells.Item(7, 4).Value = '3.633.89'
$ws.Cells.Item(7, 5).Value = '  +3.24%  '
$ws.Cells.Item(8, 4).Value = '''1.00'
$ws.Cells.Item(8, 5).Value = '  -0.09%  '
$ws.Cells.Item(9, 5).Value = '  +2.85%  '
$ws.Cells.Item(10, 5).Value = '  +7.42%  '
$ws.Cells.Item(11, 4).Value = '''7.27'
$ws.Cells.Item(11, 5).Value = '  +6.48%  '
$ws.Cells.Item(12, 5).Value = '  +4.14%  '
$ws.Cells.Item(13, 4).Value = '''0.0000232'
$ws.Cells.Item(13, 5).Value = '  +6.12%  '
$ws.Cells.Item(14, 4).Value = '''33.50'
$ws.Cells.Item(14, 5).Value = '  +6.26%  '
$ws.Cells.Item(15, 4).Value = '4.247.57'
$ws.Cells.Item(15, 5).Value = '  +2.98%  '
$ws.Cells.Item(16, 4).Value = '3.634.86'
$ws.Cells.Item(16, 5).Value = '  +3.27%  '
$ws.Cells.Item(17, 4).Value = '69.734.97'
$ws.Cells.Item(17, 5).Value = '  +4.25%  '
$ws.Cells.Item(18, 5).Value = '  +0.07%  '
$ws.Cells.Item(19, 4).Value = '''6.67'
$ws.Cells.Item(19, 5).Value = '  +6.50%  '
$ws.Cells.Item(20, 4).Value = '''16.08'
$ws.Cells.Item(20, 5).Value = '  +5.01%  '
$ws.Cells.Item(21, 4).Value = '''10.27'
$ws.Cells.Item(21, 5).Value = '  +12.03%  '
$ws.Cells.Item(22, 4).Value = '''466.19'
$ws.Cells.Item(22, 5).Value = '  +5.18%  '
$ws.Cells.Item(23, 4).Value = '''0.647'
$ws.Cells.Item(23, 5).Value = '  +2.99%  '
$ws.Cells.Item(24, 4).Value = '''78.88'
$ws.Cells.Item(24, 5).Value = '  +1.57%  '
$ws.Cells.Item(25, 4).Value = '''0.0000138'
$ws.Cells.Item(25, 5).Value = '  +12.90%  '
$ws.Cells.Item(26, 5).Value = '  +5.70%  '
$ws.Cells.Item(27, 4).Value = '3.777.20'
$ws.Cells.Item(27, 5).Value = '  +3.08%  '
$ws.Cells.Item(28, 5).Value = '  +0.03%  '
$ws.Cells.Item(29, 4).Value = '''9.33'
$ws.Cells.Item(29, 5).Value = '  +14.33%  '
$ws.Cells.Item(30, 4).Value = '''2.65'
$ws.Cells.Item(30, 5).Value = '  +5.06%  '
$ws.Cells.Item(31, 4).Value = '''1.74'
$ws.Cells.Item(31, 5).Value = '  +5.29%  '
$ws.Cells.Item(32, 5).Value = '  +13.20%  '
$ws.Cells.Item(33, 4).Value = '''6.67'
$ws.Cells.Item(33, 5).Value = '  +9.00%  '
$ws.Cells.Item(34, 4).Value = '''0.999'
$ws.Cells.Item(34, 5).Value = '  -0.02%  '
$ws.Cells.Item(35, 4).Value = '''1.99'
$ws.Cells.Item(35, 5).Value = '  +6.90%  '
$ws.Cells.Item(36, 4).Value = '''26.68'
$ws.Cells.Item(36, 5).Value = '  +3.93%  '
$ws.Cells.Item(37, 4).Value = '3.627.11'
$ws.Cells.Item(37, 5).Value = '  +3.09%  '
$ws.Cells.Item(38, 4).Value = '''8.51'
$ws.Cells.Item(38, 5).Value = '  +6.61%  '
$ws.Cells.Item(39, 4).Value = '''2.45'
$ws.Cells.Item(39, 5).Value = '  +15.75%  '
$ws.Cells.Item(40, 5).Value = '  -0.04%  '
$ws.Cells.Item(41, 4).Value = '''0.0935'
$ws.Cells.Item(41, 5).Value = '  +9.07%  '
$ws.Cells.Item(42, 2).Value = 'Monero'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(42, 4).Value = '''177.73'
$ws.Cells.Item(42, 5).Value = '  +2.04%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).Value = '''0.998'
$ws.Cells.Item(43, 5).Value = '  -0.23%  '
$ws.Cells.Item(44, 4).Value = '''5.66'
$ws.Cells.Item(44, 5).Value = '  +2.28%  '
$ws.Cells.Item(45, 4).Value = '''31.89'
$ws.Cells.Item(45, 5).Value = '  +18.15%  '
$ws.Cells.Item(46, 4).Value = '''0.918'
$ws.Cells.Item(46, 5).Value = '  +3.41%  '
$ws.Cells.Item(47, 2).Value = 'ONDO'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(47, 4).Value = '''1.38'
$ws.Cells.Item(47, 5).Value = '  +13.56%  '
$ws.Cells.Item(48, 2).Value = 'dogwifhat'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(48, 4).Value = '''2.83'
$ws.Cells.Item(48, 5).Value = '  +11.82%  '
$ws.Cells.Item(49, 4).Value = '''46.56'
$ws.Cells.Item(49, 5).Value = '  +2.99%  '
$ws.Cells.Item(50, 4).Value = '''7.83'
$ws.Cells.Item(50, 5).Value = '  +3.92%  '
$ws.Cells.Item(51, 4).Value = '''0.270'
$ws.Cells.Item(51, 5).Value = '  +9.05%  '

Write-Output "Applied cryptos update: 100 cell writes"
